# Insert a new "Skill Description" column between SkillCode (A) and SFIA Level (old B, now C)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns("B:B").Insert()

$ws.Range("B1").Value = "Skill Description"

# Generic SFIA attribute rows: Skill Description repeats the SkillCode value
$ws.Range("B2:B5").Value = "Autonomy"
$ws.Range("B6:B8").Value = "Influence"
$ws.Range("B9:B11").Value = "Complexity"
$ws.Range("B12:B14").Value = "Knowledge"

# Named skill rows: Skill Description is the human readable name of the skill code
$ws.Range("B15:B17").Value = "User research"
$ws.Range("B18:B20").Value = "User experience analysis"
$ws.Range("B21:B27").Value = "User experience evaluation"
$ws.Range("B29:B31").Value = "Methods and tools"
